# Update the workbook for release "mines - version 1.0.0 (Feb 3 2026)"
# Replaces the old build string with the new one across the "About" sheet
# (A2, A6) and the "Boundaries and methane sources" sheet (S2:S18).

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: $newVersion"

$about.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Yuwu Coal Mine, China, M0421, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

# Find last used row to safely cover the build_version column (S)
$lastRow = $data.Cells.Item($data.Rows.Count, "S").End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $data.Cells.Item($r, 19)
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}

$wb.Save()
